$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.270.16"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.380.85"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  +4.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0789"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "2.738.79"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "2.357.23"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.804"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "43.168.81"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.45%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("E34").Value = "  +10.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0734"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "128.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").Value = "1.932.47"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E46").Value = "  -8.50%  "
$ws.Range("D47").Value = "2.599.00"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
